# Update (Removed Auto Arima)
# Updates forecast values in "Forecast Comparison" sheet (C:G, rows 2-17)
# and the corresponding roll-up figures in the "Summary" sheet (B9:B14)
# now that the Auto ARIMA model has been removed from the forecast mix.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Forecast Comparison
# Columns: C = Prophet Forecast, D = Amazon Mean Forecast,
#          E = Amazon P70 Forecast, F = Amazon P80 Forecast,
#          G = Amazon P90 Forecast
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$forecastValues = @{
    2  = @(168, 40, 50, 62, 83)
    3  = @(143, 36, 44, 56, 77)
    4  = @(101, 40, 49, 62, 84)
    5  = @(70,  40, 49, 62, 83)
    6  = @(60,  41, 50, 66, 93)
    7  = @(63,  41, 50, 66, 92)
    8  = @(68,  41, 50, 65, 92)
    9  = @(70,  40, 49, 65, 91)
    10 = @(67,  39, 48, 63, 88)
    11 = @(58,  39, 48, 63, 89)
    12 = @(42,  39, 47, 63, 89)
    13 = @(28,  42, 52, 68, 95)
    14 = @(27,  39, 47, 62, 88)
    15 = @(38,  38, 46, 62, 88)
    16 = @(49,  39, 47, 62, 88)
    17 = @(48,  38, 47, 62, 88)
}

foreach ($row in $forecastValues.Keys) {
    $vals = $forecastValues[$row]
    $wsForecast.Range("C$row").Value = $vals[0]
    $wsForecast.Range("D$row").Value = $vals[1]
    $wsForecast.Range("E$row").Value = $vals[2]
    $wsForecast.Range("F$row").Value = $vals[3]
    $wsForecast.Range("G$row").Value = $vals[4]
}

# ---------------------------------------------------------------------------
# Sheet: Summary
# B9  = Total Forecast (16 Weeks)
# B10 = Total Forecast (8 Weeks)
# B11 = Total Forecast (4 Weeks)
# B12 = Max Forecast
# B14 = Min Forecast
# These cells are stored as text, so force a text number format before
# writing the value to keep the cell type consistent with the original.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$summaryValues = @{
    "B9"  = "1100"
    "B10" = "743"
    "B11" = "482"
    "B12" = "168"
    "B14" = "27"
}

foreach ($cellRef in $summaryValues.Keys) {
    $cell = $wsSummary.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryValues[$cellRef]
}
